# Updated hermite spline a little and added node class
# -> Adds sample "circle" data points (rows 11-16) and a couple of
#    trig sanity-check cells (rows 18, 19, 21) to Sheet3, then leaves
#    the selection parked just below the new data (A22), matching the
#    author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Quarter-circle sample points: x goes 0 -> 2 in steps of 0.4,
# y = sqrt(4 - x^2)  (radius 2 circle)
$ws.Range("A11").Value = 0
$ws.Range("B11").Formula = "=SQRT(4-A11^2)"

$ws.Range("A12").Formula = "=A11+0.4"
$ws.Range("B12").Formula = "=SQRT(4-A12^2)"

$ws.Range("A13").Formula = "=A12+0.4"
$ws.Range("B13").Formula = "=SQRT(4-A13^2)"

$ws.Range("A14").Formula = "=A13+0.4"
$ws.Range("B14").Formula = "=SQRT(4-A14^2)"

$ws.Range("A15").Formula = "=A14+0.4"
$ws.Range("B15").Formula = "=SQRT(4-A15^2)"

$ws.Range("A16").Formula = "=A15+0.4"
$ws.Range("B16").Formula = "=SQRT(4-A16^2)"

# A couple of trig sanity checks
$ws.Range("A18").Formula = "=2*SIN(PI()/4)"
$ws.Range("A19").Formula = "=2*COS(PI()/4)"
$ws.Range("A21").Formula = "=PI()*4/4"

# Leave the cursor just below the new block, as in the authored file.
$null = $ws.Range("A22").Select()
